# Update the "想去人数" (people interested) counts on the
# "展览" and "全部类型" sheets to match the latest scrape.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row -> new value for column F
$updates = @{
    2  = 133
    3  = 412
    4  = 11974
    5  = 1257
    11 = 407
    17 = 1790
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
